$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.94"
$ws.Range("E2").Value = "'-0.11%"
$ws.Range("E3").Value = "'-1.79%"
$ws.Range("D4").Value = "'4.720"
$ws.Range("E4").Value = "'0.01%"
$ws.Range("D5").Value = "'0.06216"
$ws.Range("E5").Value = "'2.23%"
$ws.Range("D6").Value = "'6.727"
$ws.Range("E6").Value = "'0.71%"
$ws.Range("D7").Value = "'0.8491"
$ws.Range("E7").Value = "'0.41%"
$ws.Range("D8").Value = "'0.9139"
$ws.Range("E8").Value = "'-0.82%"
$ws.Range("D9").Value = "'0.1401"
$ws.Range("E9").Value = "'-0.02%"
$ws.Range("D10").Value = "'0.04928"
$ws.Range("E10").Value = "'-2.28%"
$ws.Range("D11").Value = "'0.07091"
$ws.Range("E11").Value = "'-0.18%"
$ws.Range("D12").Value = "'0.03112"
$ws.Range("E12").Value = "'-0.61%"
$ws.Range("E13").Value = "'-0.15%"
$ws.Range("D14").Value = "'0.001535"
$ws.Range("E14").Value = "'0.01%"
$ws.Range("D15").Value = "'0.0006154"
$ws.Range("E15").Value = "'1.07%"
$ws.Range("D16").Value = "'0.005969"
$ws.Range("E16").Value = "'-2.33%"
$ws.Range("D17").Value = "'3.449"
$ws.Range("E17").Value = "'-0.09%"
$ws.Range("D18").Value = "'3.174"
$ws.Range("E18").Value = "'0.81%"
$ws.Range("D20").Value = "'0.3097"
$ws.Range("E20").Value = "'-0.94%"
$ws.Range("D21").Value = "'0.1311"
$ws.Range("E21").Value = "'0.39%"
$ws.Range("D22").Value = "'4.112"
$ws.Range("E22").Value = "'0.69%"
$ws.Range("E23").Value = "'0.52%"
$ws.Range("D24").Value = "'0.001183"
$ws.Range("E24").Value = "'-3.15%"
$ws.Range("D25").Value = "'0.004068"
$ws.Range("E25").Value = "'4.02%"
$ws.Range("E26").Value = "'0.02%"
$ws.Range("E27").Value = "'4.09%"
$ws.Range("D40").Value = "'0.03936"
$ws.Range("E40").Value = "'1.78%"
$ws.Range("D41").Value = "'0.1112"
$ws.Range("E41").Value = "'-0.11%"
$ws.Range("D42").Value = "'0.004135"
$ws.Range("E42").Value = "'0.91%"
$ws.Range("D43").Value = "'0.002142"
$ws.Range("E43").Value = "'-3.38%"
$ws.Range("D45").Value = "'0.00005163"
$ws.Range("E45").Value = "'-2.04%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D48").Value = "'0.2499"
$ws.Range("E48").Value = "'84.69%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.0002001"
